# Activities & Contacts changes - 16 Apr 2025
#
# 1. On the "Campaign" sheet, the HLSubGroup test value (cell E2) changes
#    from "CM" to "BAS".
# 2. The "Campaign" tab becomes the active/selected sheet (it was previously
#    the "Activity" tab that was selected).

$wb = $excel.ActiveWorkbook

$campaign = $wb.Worksheets.Item("Campaign")
$campaign.Range("E2").Value = "BAS"

# Make the Campaign sheet the active tab.
$campaign.Activate()
